$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "72.300.37"
$ws.Range("E2").Value = "  -0.34%  "
# Row 3
$ws.Range("D3").Value = "2.657.64"
$ws.Range("E3").Value = "  +1.12%  "
# Row 4
$ws.Range("E4").Value = "  +0.14%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.72%  "
# Row 7
$ws.Range("E7").Value = "  +0.10%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.85%  "
# Row 9
$ws.Range("D9").Value = "2.659.84"
$ws.Range("E9").Value = "  +1.34%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.27%  "
# Row 11
$ws.Range("E11").Value = "  +1.97%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.02%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.42%  "
# Row 14
$ws.Range("D14").Value = "3.147.54"
$ws.Range("E14").Value = "  +1.63%  "
# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000186"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.49%  "
# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "72.192.62"
$ws.Range("E16").Value = "  -0.18%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.24%  "
# Row 18
$ws.Range("D18").Value = "2.661.68"
$ws.Range("E18").Value = "  +1.24%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.59%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.78%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.28%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.67%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.80%  "
# Row 25
$ws.Range("E25").Value = "  +0.00%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.78%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.40%  "
# Row 28
$ws.Range("D28").Value = "2.799.00"
$ws.Range("E28").Value = "  +1.42%  "
# Row 29
$ws.Range("E29").Value = "  +0.04%  "
# Row 30
$ws.Range("D30").Value = "0.0₃0954"
$ws.Range("E30").Value = "  -0.74%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.53%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "495.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.23%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.04%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.01%  "
# Row 35
$ws.Range("E35").Value = "  +0.18%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.47%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.43%  "
# Row 38
$ws.Range("E38").Value = "  +2.29%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.05%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.89%  "
# Row 42
$ws.Range("E42").Value = "  -5.40%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.82%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.37%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.329"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.30%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "155.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.16%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.26%  "
# Row 48
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "
# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.551"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.32%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.48%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0755"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.09%  "
